$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"

$ws.Range("F2").Select()
